# Auto-generated Excel COM-interop script
# Applies the scheduled market-data refresh described by the commit diff
# (updates to currentAveragePrice*/LevePrice*/LeveProfit* columns H-N across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 392.16666
$ws.Range("I18").Value = 392.16666
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 392.16666
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -108.16666
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 9204.541999999999
$ws.Range("I33").Value = 10719.35
$ws.Range("K33").Value = 10719.35
$ws.Range("M33").Value = -10490.35
$ws.Range("H112").Value = 3438.889
$ws.Range("J112").Value = 4666.6665
$ws.Range("L112").Value = 13999.9995
$ws.Range("N112").Value = -16215.9995
$ws.Range("H131").Value = 3984.2415
$ws.Range("I131").Value = 1660.9333
$ws.Range("J131").Value = 6473.5
$ws.Range("K131").Value = 4982.7999
$ws.Range("L131").Value = 19420.5
$ws.Range("M131").Value = 57.20010000000002
$ws.Range("N131").Value = -29500.5
$ws.Range("H132").Value = 11851.895
$ws.Range("I132").Value = 11851.895
$ws.Range("K132").Value = 35555.685
$ws.Range("M132").Value = -33025.685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 41670730
$ws.Range("I122").Value = 2496.3333
$ws.Range("J122").Value = 66671668
$ws.Range("K122").Value = 7488.999899999999
$ws.Range("L122").Value = 200015004
$ws.Range("M122").Value = -5038.999899999999
$ws.Range("N122").Value = -200019904
$ws.Range("H132").Value = 27074470
$ws.Range("I132").Value = 10883.3545
$ws.Range("K132").Value = 32650.0635
$ws.Range("M132").Value = -30120.0635

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H133").Value = 100390
$ws.Range("I133").Value = 100000
$ws.Range("K133").Value = 100000
$ws.Range("M133").Value = -94940
$ws.Range("H134").Value = 2107.35
$ws.Range("I134").Value = 1788.0883
$ws.Range("K134").Value = 5364.2649
$ws.Range("M134").Value = -2829.2649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23150092
$ws.Range("I31").Value = 1834.32
$ws.Range("J31").Value = 43105490
$ws.Range("K31").Value = 1834.32
$ws.Range("L31").Value = 43105490
$ws.Range("M31").Value = -1539.32
$ws.Range("N31").Value = -43106080
$ws.Range("H34").Value = 23150092
$ws.Range("I34").Value = 1834.32
$ws.Range("J34").Value = 43105490
$ws.Range("K34").Value = 1834.32
$ws.Range("L34").Value = 43105490
$ws.Range("M34").Value = -1632.32
$ws.Range("N34").Value = -43105894
$ws.Range("H86").Value = 5612.864
$ws.Range("I86").Value = 3974.3125
$ws.Range("J86").Value = 9982.333000000001
$ws.Range("K86").Value = 3974.3125
$ws.Range("L86").Value = 9982.333000000001
$ws.Range("M86").Value = -2851.3125
$ws.Range("N86").Value = -12228.333
$ws.Range("H89").Value = 5612.864
$ws.Range("I89").Value = 3974.3125
$ws.Range("J89").Value = 9982.333000000001
$ws.Range("K89").Value = 19871.5625
$ws.Range("L89").Value = 49911.665
$ws.Range("M89").Value = -14255.5625
$ws.Range("N89").Value = -61143.665
$ws.Range("H105").Value = 10879.471
$ws.Range("I105").Value = 5129.3335
$ws.Range("J105").Value = 54005.5
$ws.Range("K105").Value = 5129.3335
$ws.Range("L105").Value = 54005.5
$ws.Range("M105").Value = -3382.3335
$ws.Range("N105").Value = -57499.5
$ws.Range("H107").Value = 503.8
$ws.Range("I107").Value = 442.3125
$ws.Range("J107").Value = 749.75
$ws.Range("K107").Value = 442.3125
$ws.Range("L107").Value = 749.75
$ws.Range("M107").Value = 1477.6875
$ws.Range("N107").Value = -4589.75
$ws.Range("H110").Value = 65000
$ws.Range("J110").Value = 65000
$ws.Range("L110").Value = 65000
$ws.Range("N110").Value = -73180
$ws.Range("H121").Value = 38974.25
$ws.Range("J121").Value = 27299
$ws.Range("L121").Value = 27299
$ws.Range("N121").Value = -29919
$ws.Range("H122").Value = 2507378.5
$ws.Range("I122").Value = 989.2353000000001
$ws.Range("J122").Value = 13159533
$ws.Range("K122").Value = 2967.7059
$ws.Range("L122").Value = 39478599
$ws.Range("M122").Value = -517.7058999999999
$ws.Range("N122").Value = -39483499
$ws.Range("H132").Value = 71950.27
$ws.Range("I132").Value = 79799.42
$ws.Range("J132").Value = 3924.3333
$ws.Range("K132").Value = 239398.26
$ws.Range("L132").Value = 11772.9999
$ws.Range("M132").Value = -236868.26
$ws.Range("N132").Value = -16832.9999
$ws.Range("H141").Value = 109475.664
$ws.Range("J141").Value = 119160.125
$ws.Range("L141").Value = 119160.125
$ws.Range("N141").Value = -129520.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3523.697
$ws.Range("J68").Value = 3899.9524
$ws.Range("L68").Value = 11699.8572
$ws.Range("N68").Value = -13321.8572
$ws.Range("H71").Value = 3523.697
$ws.Range("J71").Value = 3899.9524
$ws.Range("L71").Value = 35099.5716
$ws.Range("N71").Value = -43211.5716
$ws.Range("H132").Value = 2783088
$ws.Range("I132").Value = 3609.818
$ws.Range("K132").Value = 32488.362
$ws.Range("M132").Value = -29958.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 31665
$ws.Range("J52").Value = 31665
$ws.Range("L52").Value = 31665
$ws.Range("N52").Value = -32183
$ws.Range("H70").Value = 263000.5
$ws.Range("I70").Value = 297714.84
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 297714.84
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = -297444.84
$ws.Range("N70").Value = -20540
$ws.Range("H73").Value = 263000.5
$ws.Range("I73").Value = 297714.84
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 297714.84
$ws.Range("L73").Value = 20000
$ws.Range("M73").Value = -296778.84
$ws.Range("N73").Value = -21872
$ws.Range("H97").Value = 871.1111
$ws.Range("I97").Value = 1125.4
$ws.Range("J97").Value = 553.25
$ws.Range("K97").Value = 1125.4
$ws.Range("L97").Value = 553.25
$ws.Range("M97").Value = -629.4000000000001
$ws.Range("N97").Value = -1545.25
$ws.Range("H132").Value = 3052.7073
$ws.Range("I132").Value = 2786.7932
$ws.Range("J132").Value = 3695.3333
$ws.Range("K132").Value = 8360.3796
$ws.Range("L132").Value = 11085.9999
$ws.Range("M132").Value = -5830.3796
$ws.Range("N132").Value = -16145.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3647.36
$ws.Range("J22").Value = 4210.625
$ws.Range("L22").Value = 4210.625
$ws.Range("N22").Value = -4800.625
$ws.Range("H27").Value = 3647.36
$ws.Range("J27").Value = 4210.625
$ws.Range("L27").Value = 4210.625
$ws.Range("N27").Value = -4424.625
$ws.Range("H46").Value = 1965.1621
$ws.Range("I46").Value = 989.6786
$ws.Range("K46").Value = 989.6786
$ws.Range("M46").Value = -801.6786

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H113").Value = 1157.1538
$ws.Range("I113").Value = 1227.3
$ws.Range("J113").Value = 923.3333
$ws.Range("K113").Value = 3681.9
$ws.Range("L113").Value = 2769.9999
$ws.Range("M113").Value = -1511.9
$ws.Range("N113").Value = -7109.9999
$ws.Range("H122").Value = 28574890
$ws.Range("I122").Value = 31252756
$ws.Range("J122").Value = 20005722
$ws.Range("K122").Value = 93758268
$ws.Range("L122").Value = 60017166
$ws.Range("M122").Value = -93755818
$ws.Range("N122").Value = -60022066

